$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (pop1 / Automation_1 first data row) updates ---
$ws.Range("C2").Value = "NewImportLogic_2 - Test_Automation_2"
$ws.Range("D2").Value = "Success Sheet - LineofTherapy_1.xlsx"
$ws.Range("E2").Value = "\Testdata\Templates\LineOfTherapy\Testing_Env\Success Sheet - LineofTherapy_1.xlsx"
$ws.Range("F2").Value = 48
$ws.Range("I2").Value = "Manage Population filter 2"

# --- Row 3: new pop1 error-row data ---
$ws.Range("A3").Value = "pop1"
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = "Population filter 2 'Automation_1' is not supported"
$ws.Range("I3").Value = "You can view all, create new and edit or delete existing Population filter 2 from here"

# --- Row 4: now pop1 (was pop2); B4/C4/D4/E4 cleared ---
$ws.Range("A4").Value = "pop1"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 122
$ws.Range("G4").Value = "Population filter 2 'Automation_1' is not supported"

# --- Row 5: new pop1 error-row data ---
$ws.Range("A5").Value = "pop1"
$ws.Range("F5").Value = 146
$ws.Range("G5").Value = "Population filter 2 'Automation_1' is not supported"

# --- Row 7: new pop2 / Automation_2 data row ---
$ws.Range("A7").Value = "pop2"
$ws.Range("B7").Value = "Automation_2"
$ws.Range("C7").Value = "NewImportLogic_2 - Test_Automation_2"
$ws.Range("D7").Value = "Success Sheet - LineofTherapy_2.xlsx"
$ws.Range("E7").Value = "\Testdata\Templates\LineOfTherapy\Testing_Env\Success Sheet - LineofTherapy_2.xlsx"
$ws.Range("F7").Value = 28
$ws.Range("G7").Value = "Population filter 2 'Automation_2' is not supported"

# --- Row 8: new pop2 error-row data ---
$ws.Range("A8").Value = "pop2"
$ws.Range("F8").Value = 62
$ws.Range("G8").Value = "Population filter 2 'Automation_2' is not supported"

# --- Row 9: new pop2 error-row data ---
$ws.Range("A9").Value = "pop2"
$ws.Range("F9").Value = 114
$ws.Range("G9").Value = "Population filter 2 'Automation_2' is not supported"

# --- Row 10: brand-new row ---
$ws.Range("A10").Value = "pop2"
$ws.Range("F10").Value = 152
$ws.Range("G10").Value = "Population filter 2 'Automation_2' is not supported"

# --- View: selection covers A1:I11, scrolled so column E is the leftmost visible ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("A1:I11").Select()
